$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "57.439.62"
$ws.Range("E2").Value2 = "  -4.79%  "

$ws.Range("D3").Value2 = "3.119.67"
$ws.Range("E3").Value2 = "  -5.79%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "0.999"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "521.37"
$ws.Range("E5").Value2 = "  -6.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "134.53"
$ws.Range("E6").Value2 = "  -5.47%  "

$ws.Range("E7").Value2 = "  -0.01%  "

$ws.Range("D8").Value2 = "3.114.58"
$ws.Range("E8").Value2 = "  -5.94%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.445"
$ws.Range("E9").Value2 = "  -6.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "7.19"
$ws.Range("E10").Value2 = "  -8.42%  "

$ws.Range("E11").Value2 = "  -8.59%  "

$ws.Range("E12").Value2 = "  -6.22%  "

$ws.Range("D13").Value2 = "3.653.09"
$ws.Range("E13").Value2 = "  -5.87%  "

$ws.Range("E14").Value2 = "  -2.30%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "25.32"
$ws.Range("E15").Value2 = "  -6.11%  "

$ws.Range("D16").Value2 = "3.118.72"
$ws.Range("E16").Value2 = "  -5.67%  "

$ws.Range("D17").Value2 = "57.363.14"
$ws.Range("E17").Value2 = "  -4.89%  "

$ws.Range("E18").Value2 = "  -9.77%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "5.75"
$ws.Range("E19").Value2 = "  -6.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "12.94"
$ws.Range("E20").Value2 = "  -10.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "7.95"
$ws.Range("E21").Value2 = "  -8.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "341.65"
$ws.Range("E22").Value2 = "  -9.08%  "

$ws.Range("E23").Value2 = "  -0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "68.05"
$ws.Range("E24").Value2 = "  -8.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "0.502"
$ws.Range("E25").Value2 = "  -7.72%  "

$ws.Range("D26").Value2 = "3.247.44"
$ws.Range("E26").Value2 = "  -5.65%  "

$ws.Range("E27").Value2 = "  -3.64%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "0.998"
$ws.Range("E28").Value2 = "  -0.06%  "

$ws.Range("D29").Value2 = "0.0₃0932"
$ws.Range("E29").Value2 = "  -9.57%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "0.998"
$ws.Range("E30").Value2 = "  -0.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "6.72"
$ws.Range("E31").Value2 = "  -7.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "6.95"
$ws.Range("E32").Value2 = "  -9.63%  "

$ws.Range("E33").Value2 = "  -9.25%  "

$ws.Range("E34").Value2 = "  -3.34%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "21.39"
$ws.Range("E35").Value2 = "  -5.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "158.69"
$ws.Range("E36").Value2 = "  -4.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "4.75"
$ws.Range("E37").Value2 = "  -8.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "6.15"
$ws.Range("E38").Value2 = "  -9.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "1.37"
$ws.Range("E39").Value2 = "  -10.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "25.06"
$ws.Range("E40").Value2 = "  -6.89%  "

$ws.Range("E41").Value2 = "  -7.22%  "

$ws.Range("D42").Value2 = "3.146.21"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "40.22"
$ws.Range("E43").Value2 = "  -4.30%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.679"
$ws.Range("E44").Value2 = "  -9.72%  "

$ws.Range("B45").Value2 = "Filecoin"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "3.90"
$ws.Range("E45").Value2 = "  -7.46%  "

$ws.Range("B46").Value2 = "ONDO"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "1.06"
$ws.Range("E46").Value2 = "  -5.11%  "

$ws.Range("E47").Value2 = "  -0.07%  "

$ws.Range("E48").Value2 = "  -9.81%  "

$ws.Range("D49").Value2 = "2.254.75"
$ws.Range("E49").Value2 = "  -5.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "6.15"
$ws.Range("E50").Value2 = "  -5.98%  "

$ws.Range("E51").Value2 = "  -7.26%  "
